# Applies the Sat Oct 21 14:17:46 UTC 2023 cryptos list refresh (prices, 1h volume %, and
# a few re-ranked rows) to the active worksheet, cell by cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.683.42"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "1.607.80"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("D4").Value = "'" + "0.996"
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").Value = "'" + "212.17"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").Value = "'" + "0.518"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").Value = "'" + "0.995"
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("D8").Value = "'" + "29.00"
$ws.Range("E8").Value = "  +8.24%  "
$ws.Range("E9").Value = "  +3.62%  "
$ws.Range("D10").Value = "'" + "0.0607"
$ws.Range("E10").Value = "  +1.55%  "
$ws.Range("D11").Value = "'" + "0.0907"
$ws.Range("E11").Value = "  -0.81%  "
$ws.Range("D12").Value = "1.838.03"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").Value = "1.604.65"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("D14").Value = "'" + "0.564"
$ws.Range("E14").Value = "  +5.58%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'" + "3.83"
$ws.Range("E15").Value = "  +2.59%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "29.686.68"
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("D17").Value = "'" + "64.56"
$ws.Range("E17").Value = "  +1.77%  "
$ws.Range("D18").Value = "'" + "8.36"
$ws.Range("E18").Value = "  +9.87%  "
$ws.Range("D19").Value = "'" + "241.39"
$ws.Range("E19").Value = "  +0.91%  "
$ws.Range("D20").Value = "0.0₃0702"
$ws.Range("E20").Value = "  +1.49%  "
$ws.Range("D21").Value = "'" + "0.996"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("E22").Value = "  +1.00%  "
$ws.Range("E23").Value = "  +4.02%  "
$ws.Range("D24").Value = "'" + "2.11"
$ws.Range("E24").Value = "  +1.77%  "
$ws.Range("D25").Value = "'" + "156.62"
$ws.Range("E25").Value = "  +1.29%  "
$ws.Range("D26").Value = "'" + "15.54"
$ws.Range("E26").Value = "  +1.78%  "
$ws.Range("E27").Value = "  +1.13%  "
$ws.Range("D28").Value = "'" + "6.53"
$ws.Range("E28").Value = "  +2.60%  "
$ws.Range("D29").Value = "'" + "0.997"
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("E30").Value = "  +2.03%  "
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("D32").Value = "'" + "3.27"
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("D33").Value = "'" + "3.18"
$ws.Range("E33").Value = "  +2.43%  "
$ws.Range("D34").Value = "1.424.05"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("E35").Value = "  +5.07%  "
$ws.Range("E36").Value = "  +1.06%  "
$ws.Range("D37").Value = "'" + "2.88"
$ws.Range("E37").Value = "  +2.10%  "
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("E39").Value = "  +2.73%  "
$ws.Range("D40").Value = "'" + "0.554"
$ws.Range("E40").Value = "  +3.74%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'" + "0.0497"
$ws.Range("E41").Value = "  +5.33%  "
$ws.Range("E42").Value = "  +1.01%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'" + "0.824"
$ws.Range("E43").Value = "  +4.01%  "
$ws.Range("D44").Value = "'" + "54.50"
$ws.Range("E44").Value = "  +2.96%  "
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "'" + "0.995"
$ws.Range("E45").Value = "  -0.26%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'" + "68.19"
$ws.Range("E46").Value = "  +4.43%  "
$ws.Range("D47").Value = "'" + "0.997"
$ws.Range("E47").Value = "  +19.26%  "
$ws.Range("D48").Value = "'" + "5.42"
$ws.Range("E48").Value = "  +2.64%  "
$ws.Range("D49").Value = "1.746.97"
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("D50").Value = "'" + "87.21"
$ws.Range("E50").Value = "  +0.89%  "
$ws.Range("E51").Value = "  -1.67%  "
